$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(1)
$p.Range.Font.Bold = $p.Range.Font.Bold
